$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a number by Excel
# are force-formatted as Text, written, then restored to General so the
# stored cell keeps looking like a plain (unformatted) text cell,
# matching the original inlineStr text values in the workbook.

$ws.Range("D2").Value = '58.846.57'
$ws.Range("E2").Value = '  +2.14%  '
$ws.Range("D3").Value = '2.588.28'
$ws.Range("E3").Value = '  +0.74%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '520.61'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = '  -0.06%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '139.91'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = '  -2.95%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("E8").Value = '  +0.62%  '
$ws.Range("D9").Value = '2.597.95'
$ws.Range("E9").Value = '  +0.56%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.53'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = '  -0.69%  '
$ws.Range("E11").Value = '  +0.28%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.331'
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = '  +1.54%  '
$ws.Range("E13").Value = '  +2.94%  '
$ws.Range("D14").Value = '3.044.68'
$ws.Range("E14").Value = '  +0.69%  '
$ws.Range("D15").Value = '58.796.56'
$ws.Range("E15").Value = '  +2.07%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '20.48'
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = '  +1.33%  '
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '2.593.19'
$ws.Range("E17").Value = '  +0.86%  '
$ws.Range("B18").Value = 'ShibaInu'
$ws.Range("C18").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0000133'
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = '  -0.24%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '338.69'
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = '  +0.88%  '
$ws.Range("E20").Value = '  +0.31%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.19'
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = '  +0.10%  '
$ws.Range("E22").Value = '  +3.45%  '
$ws.Range("E23").Value = '  -0.03%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '66.10'
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = '  +2.20%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.168'
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = '  +0.65%  '
$ws.Range("E26").Value = '  +0.91%  '
$ws.Range("E27").Value = '  +0.16%  '
$ws.Range("E28").Value = '  +1.27%  '
$ws.Range("E29").Value = '  +0.06%  '
$ws.Range("D30").Value = '0.0₃0724'
$ws.Range("E30").Value = '  -3.53%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.93'
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = '  -5.28%  '
$ws.Range("E32").Value = '  -1.11%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '18.77'
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = '  +0.85%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '148.99'
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = '  +0.41%  '
$ws.Range("E35").Value = '  -1.42%  '
$ws.Range("E36").Value = '  -1.51%  '
$ws.Range("B37").Value = 'OKB'
$ws.Range("C37").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '36.29'
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = '  +1.11%  '
$ws.Range("B38").Value = 'Stacks'
$ws.Range("C38").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.47'
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = '  +1.97%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.828'
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = '  -1.15%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.820'
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = '  -2.47%  '
$ws.Range("E41").Value = '  -0.08%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.997'
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = '  -0.05%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '274.31'
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = '  +1.87%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '10.74'
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = '  +0.94%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.589'
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = '  -0.10%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0951'
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = '  -0.30%  '
$ws.Range("E47").Value = '  +0.04%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '18.55'
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = '  -1.54%  '
$ws.Range("D49").Value = '1.990.47'
$ws.Range("E49").Value = '  +0.83%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0220'
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = '  +0.20%  '
$ws.Range("E51").Value = '  -4.75%  '
